$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2028571428571428
$ws.Range("C2").Value = 0.5428571428571428
$ws.Range("J2").Value = 0.02571428571428571
$ws.Range("P2").Value = 0.1285714285714286
$ws.Range("S2").Value = 0.1
$ws.Range("C3").Value = 0.03553299492385787
$ws.Range("J3").Value = 0.04060913705583756
$ws.Range("P3").Value = 0.700507614213198
$ws.Range("S3").Value = 0.2233502538071066
$ws.Range("J4").Value = 0.1
$ws.Range("P4").Value = 0.6333333333333333
$ws.Range("S4").Value = 0.2666666666666667
$ws.Range("B6").Value = 0.0958904109589041
$ws.Range("D6").Value = 0.0136986301369863
$ws.Range("E6").Value = 0.0045662100456621
$ws.Range("F6").Value = 0.0547945205479452
$ws.Range("J6").Value = 0.2648401826484018
$ws.Range("O6").Value = 0.0228310502283105
$ws.Range("Q6").Value = 0.1643835616438356
$ws.Range("R6").Value = 0.045662100456621
$ws.Range("S6").Value = 0.3333333333333333
$ws.Range("B7").Value = 0.1095890410958904
$ws.Range("D7").Value = 0.045662100456621
$ws.Range("E7").Value = 0.0045662100456621
$ws.Range("F7").Value = 0.0502283105022831
$ws.Range("J7").Value = 0.1506849315068493
$ws.Range("O7").Value = 0.0136986301369863
$ws.Range("Q7").Value = 0.2054794520547945
$ws.Range("R7").Value = 0.1050228310502283
$ws.Range("S7").Value = 0.3150684931506849
$ws.Range("B8").Value = 0.09848484848484848
$ws.Range("D8").Value = 0.0303030303030303
$ws.Range("E8").Value = 0.002525252525252525
$ws.Range("F8").Value = 0.07323232323232323
$ws.Range("J8").Value = 0.154040404040404
$ws.Range("O8").Value = 0.01262626262626263
$ws.Range("Q8").Value = 0.1666666666666667
$ws.Range("R8").Value = 0.0707070707070707
$ws.Range("S8").Value = 0.3914141414141414
$ws.Range("B9").Value = 0.1386138613861386
$ws.Range("D9").Value = 0.009900990099009901
$ws.Range("E9").Value = 0.004950495049504951
$ws.Range("F9").Value = 0.05445544554455446
$ws.Range("J9").Value = 0.1584158415841584
$ws.Range("O9").Value = 0.03465346534653466
$ws.Range("Q9").Value = 0.1336633663366337
$ws.Range("R9").Value = 0.0594059405940594
$ws.Range("S9").Value = 0.405940594059406
$ws.Range("B10").Value = 0.1147887323943662
$ws.Range("D10").Value = 0.02464788732394366
$ws.Range("E10").Value = 0.0007042253521126761
$ws.Range("F10").Value = 0.06549295774647887
$ws.Range("J10").Value = 0.1450704225352113
$ws.Range("O10").Value = 0.01197183098591549
$ws.Range("Q10").Value = 0.2394366197183098
$ws.Range("R10").Value = 0.07605633802816901
$ws.Range("S10").Value = 0.321830985915493
$ws.Range("G11").Value = 0.1552238805970149
$ws.Range("J11").Value = 0.08059701492537313
$ws.Range("K11").Value = 0.182089552238806
$ws.Range("L11").Value = 0.564179104477612
$ws.Range("S11").Value = 0.01791044776119403
$ws.Range("G12").Value = 0.7135416666666666
$ws.Range("J12").Value = 0.2291666666666667
$ws.Range("K12").Value = 0.005208333333333333
$ws.Range("L12").Value = 0.03125
$ws.Range("S12").Value = 0.02083333333333333
$ws.Range("G13").Value = 0.68
$ws.Range("J13").Value = 0.28
$ws.Range("S13").Value = 0.04
$ws.Range("F15").Value = 0.03070175438596491
$ws.Range("H15").Value = 0.131578947368421
$ws.Range("I15").Value = 0.07017543859649122
$ws.Range("J15").Value = 0.3903508771929824
$ws.Range("K15").Value = 0.06578947368421052
$ws.Range("O15").Value = 0.09210526315789473
$ws.Range("S15").Value = 0.2192982456140351
$ws.Range("F16").Value = 0.009049773755656109
$ws.Range("H16").Value = 0.1719457013574661
$ws.Range("I16").Value = 0.06787330316742081
$ws.Range("J16").Value = 0.416289592760181
$ws.Range("K16").Value = 0.1131221719457014
$ws.Range("M16").Value = 0.02714932126696833
$ws.Range("O16").Value = 0.04977375565610859
$ws.Range("S16").Value = 0.1447963800904978
$ws.Range("F17").Value = 0.01171875
$ws.Range("H17").Value = 0.169921875
$ws.Range("I17").Value = 0.091796875
$ws.Range("J17").Value = 0.40234375
$ws.Range("K17").Value = 0.109375
$ws.Range("M17").Value = 0.01953125
$ws.Range("O17").Value = 0.064453125
$ws.Range("S17").Value = 0.130859375
$ws.Range("F18").Value = 0.01666666666666667
$ws.Range("H18").Value = 0.1555555555555556
$ws.Range("I18").Value = 0.1333333333333333
$ws.Range("J18").Value = 0.4166666666666667
$ws.Range("K18").Value = 0.09444444444444444
$ws.Range("M18").Value = 0.02222222222222222
$ws.Range("N18").Value = 0.005555555555555556
$ws.Range("O18").Value = 0.06666666666666667
$ws.Range("S18").Value = 0.08888888888888889
$ws.Range("F19").Value = 0.01610305958132045
$ws.Range("H19").Value = 0.1731078904991948
$ws.Range("I19").Value = 0.07890499194847021
$ws.Range("J19").Value = 0.3848631239935588
$ws.Range("K19").Value = 0.1296296296296296
$ws.Range("M19").Value = 0.02657004830917874
$ws.Range("N19").Value = 0.002415458937198068
$ws.Range("O19").Value = 0.06682769726247988
$ws.Range("S19").Value = 0.1215780998389694
